# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G ("K") rows 2-18
$gValues = @{
    2  = 0
    3  = 1
    4  = 0
    5  = 2
    6  = 0
    7  = 1
    8  = 2
    9  = 0
    10 = 2
    11 = 0
    12 = 0
    13 = 0
    14 = 2
    15 = 0
    16 = 2
    17 = 0
    18 = 4
}

foreach ($row in $gValues.Keys) {
    $ws.Range("G$row").Value = $gValues[$row]
}
